# Weekly update: a new week of price data is published for
# "Pepino ensalada" (Agrícola del Norte S.A. de Arica).
#
# The consolidated sheet keeps its history ordered from most-recent to
# oldest, so the new week's two rows (Primera / Segunda quality) are
# inserted at the top of the data block (row 147), pushing every
# existing data row down by two positions. This naturally reproduces
# the two extra rows at the bottom of the sheet (the oldest week,
# previously on rows 264:265, now lives on rows 266:267) without any
# further action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new week by inserting two rows above the first
# data row of the existing block (row 147). Formatting (date style on
# column D, etc.) is inherited from the row below, exactly like doing
# this by hand in Excel.
$ws.Rows("147:148").Insert()

# Row 147: Pepino ensalada, Sin especificar, Primera
$ws.Cells.Item(147, 1).Value = 1
$ws.Cells.Item(147, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(147, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(147, 4).Value = 44634
$ws.Cells.Item(147, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(147, 5).Value = 15
$ws.Cells.Item(147, 6).Value = 100112043
$ws.Cells.Item(147, 7).Value = "Pepino ensalada"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 120
$ws.Cells.Item(147, 11).Value = 14000
$ws.Cells.Item(147, 12).Value = 15000
$ws.Cells.Item(147, 13).Value = 14500
$ws.Cells.Item(147, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(147, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(147, 16).Value = 207
$ws.Cells.Item(147, 17).Value = 70
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Row 148: Pepino ensalada, Sin especificar, Segunda
$ws.Cells.Item(148, 1).Value = 1
$ws.Cells.Item(148, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(148, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(148, 4).Value = 44634
$ws.Cells.Item(148, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(148, 5).Value = 15
$ws.Cells.Item(148, 6).Value = 100112043
$ws.Cells.Item(148, 7).Value = "Pepino ensalada"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Segunda"
$ws.Cells.Item(148, 10).Value = 130
$ws.Cells.Item(148, 11).Value = 11000
$ws.Cells.Item(148, 12).Value = 12000
$ws.Cells.Item(148, 13).Value = 11500
$ws.Cells.Item(148, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 115
$ws.Cells.Item(148, 17).Value = 100
$ws.Cells.Item(148, 18).Value = "Hortaliza"
